# Weekly update of the "Zapallo" (Hortofrutícola Agro Chillán) dataset:
# a new daily record is inserted at row 148 (everything below shifts down
# by one row), and the sheet's used range grows from A1:R230 to A1:R231.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 148; this also pushes the
# dimension/used-range down by one row automatically.
$ws.Rows.Item(148).Insert()

# Populate the new row 148 with the new record. Most fields repeat the
# values of the (now shifted-down) row immediately below it; only the
# date and the price/volume columns differ.
$ws.Cells.Item(148, 1).Value = 7
$ws.Cells.Item(148, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(148, 3).Value = "Ñuble"
$ws.Cells.Item(148, 4).Value = 44960
$ws.Cells.Item(148, 5).Value = 16
$ws.Cells.Item(148, 6).Value = 100112045
$ws.Cells.Item(148, 7).Value = "Zapallo"
$ws.Cells.Item(148, 8).Value = "Camote"
$ws.Cells.Item(148, 9).Value = "1a (cosecha)"
$ws.Cells.Item(148, 10).Value = 200
$ws.Cells.Item(148, 11).Value = 500
$ws.Cells.Item(148, 12).Value = 500
$ws.Cells.Item(148, 13).Value = 500
$ws.Cells.Item(148, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(148, 15).Value = "Región del Maule"
$ws.Cells.Item(148, 16).Value = 500
$ws.Cells.Item(148, 17).Value = 1
$ws.Cells.Item(148, 18).Value = "Hortaliza"
